# SDWe注意事项.xlsx — add a new "middle weight" sample note below the
# existing two notes: row 41 gets a new 24pt red 等线 note (a smaller
# sibling of the existing 28pt red notes used for rows 33/36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note text (becomes a new shared-string entry).
$ws.Range("A41").Value = "如果出现下载时CRC错误，可以删除VT_SET下的CheckSum.bin"

# Clone the formatting of the existing "sample" row (A33: red/28pt/等线)
# onto the new cell, then dial the size down to the "middle weight" 24pt
# so only one new font + one new cell style get created.
$ws.Range("A33").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").Font.Size = 24

# Match the taller row height used for the other note rows.
$ws.Rows.Item(41).RowHeight = 30

# Leave the selection on the newly added cell.
$ws.Range("A41").Select() | Out-Null
